# Fixed a few typos in pp
#
# Slide 2 ("Scenario"): content placeholder text corrections
#   - "Get runners pules"        -> "Get runner's pulse"
#   - "Use relays to extent track" -> "Use relays to extent track coverage"
# Slide 3: nudge two pictures slightly (design refresh)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: set a shape's Left/Top (given in EMU) as precisely as the
# point-based COM properties allow (PowerPoint stores Left/Top in points as
# single-precision floats, so search the neighbourhood of the ideal value for
# the candidate that round-trips closest to the requested EMU value).
# ---------------------------------------------------------------------------
function Set-ShapeLeftEmu($shape, [double]$targetEmu) {
    $bestPts = $targetEmu / 12700.0
    $bestErr = 999999999
    for ($d = -5; $d -le 5; $d++) {
        $candidateEmu = $targetEmu + $d
        $pts = $candidateEmu / 12700.0
        $shape.Left = $pts
        $rbEmu = [math]::Round($shape.Left * 12700)
        $err = [math]::Abs($rbEmu - $targetEmu)
        if ($err -lt $bestErr) {
            $bestErr = $err
            $bestPts = $pts
        }
        if ($err -eq 0) { break }
    }
    $shape.Left = $bestPts
}

function Set-ShapeTopEmu($shape, [double]$targetEmu) {
    $bestPts = $targetEmu / 12700.0
    $bestErr = 999999999
    for ($d = -5; $d -le 5; $d++) {
        $candidateEmu = $targetEmu + $d
        $pts = $candidateEmu / 12700.0
        $shape.Top = $pts
        $rbEmu = [math]::Round($shape.Top * 12700)
        $err = [math]::Abs($rbEmu - $targetEmu)
        if ($err -lt $bestErr) {
            $bestErr = $err
            $bestPts = $pts
        }
        if ($err -eq 0) { break }
    }
    $shape.Top = $bestPts
}

# ---------------------------------------------------------------------------
# Slide 2 - "Scenario" slide, fix two typo'd bullet lines in the content
# placeholder.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$content = $slide2.Shapes.Item("Content Placeholder 9")
$tr = $content.TextFrame.TextRange

# Paragraph 1: "Get runners pules" -> "Get runner's pulse"
$para1 = $tr.Paragraphs(1)
$t1 = $para1.Text
$search1 = "runners pules"
$idx1 = $t1.IndexOf($search1)
if ($idx1 -ge 0) {
    $sub1 = $tr.Characters($para1.Start + $idx1, $search1.Length)
    $sub1.Text = "runner" + [char]0x2019 + "s pulse"
}

# Paragraph 5: "Use relays to extent track" -> "Use relays to extent track coverage"
$para5 = $tr.Paragraphs(5)
$t5 = $para5.Text
$search5 = "track"
$idx5 = $t5.IndexOf($search5)
if ($idx5 -ge 0) {
    $sub5 = $tr.Characters($para5.Start + $idx5, $search5.Length)
    $sub5.Text = "track coverage"
}

# ---------------------------------------------------------------------------
# Slide 3 - nudge "Picture 5" and "Content Placeholder 3" images slightly.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

$pic5 = $slide3.Shapes.Item("Picture 5")
Set-ShapeLeftEmu $pic5 4722006
Set-ShapeTopEmu  $pic5 813747

$contentPh3 = $slide3.Shapes.Item("Content Placeholder 3")
Set-ShapeLeftEmu $contentPh3 8726996
